# Add data sourced from MySQL (Touroperator, Trips Sold, Income) next to the
# existing SQLite-sourced columns (Trip Name, Country), per commit:
# "Add writhing into Excel report from MySql."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header titles (row 2) - rename existing + add the 3 new MySQL columns
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Trip Name (from SQLite)"
$ws.Range("B2").Value = "Country (from SQLite)"
$ws.Range("C2").Value = "Touroperator (from MySQL)"
$ws.Range("D2").Value = "Trips Sold (from MySQL)"
$ws.Range("E2").Value = "Income (from MySQL)"

# ---------------------------------------------------------------------------
# 2. Data rows
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Italy The One And Only"
$ws.Range("B3").Value = "Europe"
$ws.Range("C3").Value = "Eleganca Tours"
$ws.Range("D3").Value = 14
$ws.Range("E3").Value = 12672

$ws.Range("A4").Value = "The Secret Beauty Of Mexico"
$ws.Range("B4").Value = "North America"
$ws.Range("C4").Value = "Elite Travel Agency"
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 40920

$ws.Range("A5").Value = "Chiloe Chile"
$ws.Range("B5").Value = "South America"
$ws.Range("C5").Value = "Elite Travel Agency"
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 25800

# The old row 6 ("Mie Perfecture Japan" / "Asia") is no longer present in the
# report, drop it entirely so the used range shrinks back down to 5 rows.
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3. Title row / merge now spans A:E instead of just A:B
# ---------------------------------------------------------------------------
$ws.Range("A1:E1").Merge()

# ---------------------------------------------------------------------------
# 4. Column widths - match column A/B's width for the 3 new columns
# ---------------------------------------------------------------------------
$ws.Range("C1:E1").ColumnWidth = 24.75

# ---------------------------------------------------------------------------
# 5. Formatting
#    - Row 1 (title) keeps bold for A1; B1:E1 become plain (not bold).
#    - Row 2 (headers) keeps bold font + teal fill across A2:E2.
#    - Whole used range gets centered horizontal alignment (applied last so
#      it is folded into each cell's final resolved style).
# ---------------------------------------------------------------------------
$ws.Range("B1:E1").Font.Bold = $false

$ws.Range("C2:E2").Font.Bold = $true
$ws.Range("C2:E2").Interior.Color = 10526303     # matches existing fill (FF5F9EA0)

$ws.Range("A1:E5").HorizontalAlignment = -4108   # xlCenter

Write-Output "done"
